# Add a new "Sheet2" after "Sheet1" containing a year -> photo-count table,
# matching the commit "added photo count to sheets".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 and make it the active sheet/tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row (bold).
$ws2.Range("A1").Value = "ANO"
$ws2.Range("B1").Value = "PHOTOS"
$ws2.Range("A1:B1").Font.Bold = $true

# Year / photo-count data.
$data = @(
    @(2011, 16),
    @(2014, 8318),
    @(2015, 27009),
    @(2016, 33789),
    @(2017, 3938),
    @(2018, 96022),
    @(2019, 512519),
    @(2020, 8539),
    @(2021, 1741455),
    @(2022, 20),
    @(2023, 1413),
    @(2024, 144)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Match the authored selection state on the new sheet.
$null = $ws2.Range("E34").Select()
